$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Append two new log rows (146 and 147) after the existing data (which ends at row 145)
$ws.Range("A146").Value = 145
$ws.Range("B146").Value = 1
$ws.Range("C146").Value = "2024-06-18 00:57:57"
$ws.Range("D146").Value = 200
$ws.Range("E146").Value = 10

$ws.Range("A147").Value = 146
$ws.Range("B147").Value = 2
$ws.Range("C147").Value = "2024-06-18 00:57:57"
$ws.Range("D147").Value = 200
$ws.Range("E147").Value = 2
